# Insert a new weekly price record for Limón (Vega Monumental Concepción)
# above the current row 176, shifting all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 176 (pushes old rows 176-277 to 177-278)
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new weekly record
$ws.Range("A176").Value = 11
$ws.Range("B176").Value = "Vega Monumental Concepción"
$ws.Range("C176").Value = "Bíobío"
$ws.Range("D176").Value = 44488
$ws.Range("E176").Value = 8
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100102
$ws.Range("H176").Value = "Cítricos"
$ws.Range("I176").Value = 100102003
$ws.Range("J176").Value = "Limón"
$ws.Range("K176").Value = "Sin especificar"
$ws.Range("L176").Value = "1a amarillo"
$ws.Range("M176").Value = 430
$ws.Range("N176").Value = 6000
$ws.Range("O176").Value = 6500
$ws.Range("P176").Value = 6233
$ws.Range("Q176").Value = "$/malla 18 kilos"
$ws.Range("R176").Value = "Región de O'Higgins"
$ws.Range("S176").Value = 346
$ws.Range("T176").Value = 18
